$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update schemeMap column (M) for rows 2-4 from numeric 2 to text "HTTPS"
$ws.Range("M2").Value = "HTTPS"
$ws.Range("M3").Value = "HTTPS"
$ws.Range("M4").Value = "HTTPS"

# Move the selection to M4, matching the new active cell in the source file
$ws.Range("M4").Select()
